$d = $word.ActiveDocument

# 1. Fix hyphenation: "bem sucedidos" -> "bem-sucedidos"
$d.Content.Find.Execute(
    "querem ser bem sucedidos na sua profissão",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "querem ser bem-sucedidos na sua profissão",
    2)

# 2. Add trailing space after "dessa empresa."
$d.Content.Find.Execute(
    "dos objetivos de negócio dessa empresa.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "dos objetivos de negócio dessa empresa. ",
    2)

# 3. Add trailing space after "inovação de produto."
$d.Content.Find.Execute(
    "necessárias para a inovação de produto.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "necessárias para a inovação de produto. ",
    2)

# 4. Add trailing space after "que você é um"
$d.Content.Find.Execute(
    "Digamos, por exemplo, que você é um",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Digamos, por exemplo, que você é um ",
    2)
